$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "3Sum" row of data (row 16)
$ws.Range("A16").Value = "3Sum"
$ws.Range("B16").Value = "Return triplet of sum 0"
$ws.Range("C16").Value = "Sort input array. Use for loop to go over elements. Only consider negatives and 0 as a potential first element. Use left/right pointer to find other two elements same as Two Sum II."
$ws.Range("D16").Value = "https://leetcode.com/problems/3sum"

# Hyperlink the new leetcode URL cell (mirrors the other Link-column cells)
$ws.Hyperlinks.Add($ws.Range("D16"), "https://leetcode.com/problems/3sum")

# Match the formatting used by the other data rows (row 11 is a good template:
# A = "Neutral" style, B/C = body text style, D = hyperlink style)
$ws.Range("A11").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D11").Copy()
$ws.Range("D16").PasteSpecial(-4122)

# Match the final selected cell recorded in the saved workbook
$ws.Range("C8").Select() | Out-Null
